$d = $word.ActiveDocument

# The Pearson logo inline pictures (in the default/first-page footers) and
# the BTec logo inline picture (in the first-page header) are renamed:
#   - Pearson logo: image1.png -> image2.png
#   - BTec logo:    image2.jpg -> image1.jpg
# Word keeps the picture's "name" in two places per <w:drawing> -
# <wp:docPr name="..."/> and the nested <pic:cNvPr name="..."/> - and both
# need to change together, which isn't reachable through the InlineShape
# object model (InlineShape.Name only ever touches <wp:docPr>), so edit the
# underlying OOXML for the whole document directly.

$xml = $d.Content.WordOpenXML

$xml = $xml.Replace('name="image1.png"', 'name="image2.png"')
$xml = $xml.Replace('name="image2.jpg"', 'name="image1.jpg"')

$d.Content.WordOpenXML = $xml
